{"js": "// --- 1. Rewrite the \"Skilled in the ...\" summary paragraph ---------------\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  p.load(\"text\");\n}\nawait context.sync();\n\nconst summary = paragraphs.items.find((p) => p.text.indexOf(\"Skilled in the\") !== -1);\n\nconst newSummary =\n  \"Self motivated computer science student skilled in building and maintaining mobile applications while meeting time constraints. Well-versed in various algorithm design paradigms using discrete mathematics. Experienced in Object-Oriented design, analysis, and prototyping in teams of 3-6 developers.\";\nsummary.insertText(newSummary, \"Replace\");\nawait context.sync();\n\n// --- 2. Remove the trailing \"Excel Forecasting Simulations\" table row ----\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.rows.load(\"items\");\nawait context.sync();\n\nconst lastRow = table.rows.items[table.rows.items.length - 1];\nlastRow.delete();\nawait context.sync();\n\n// --- 3. Re-even the remaining column widths (3116 -> 3117 twips) ---------\ntable.rows.load(\"items\");\nawait context.sync();\n\nconst firstRow = table.rows.items[0];\nfirstRow.cells.load(\"items\");\nawait context.sync();\n\nfirstRow.cells.items[0].columnWidth = 155.85;\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- 1. Rewrite the \"Skilled in the ...\" summary paragraph ---------------\n# Collapse the whole paragraph (several runs) into the new wording while\n# keeping the paragraph mark (and the _GoBack bookmark) untouched so the\n# run-level formatting (Times New Roman / black / en-CA) carries over.\n$summary = $d.Paragraphs.Item(4)\n$newSummary = \"Self motivated computer science student skilled in building and maintaining mobile applications while meeting time constraints. Well-versed in various algorithm design paradigms using discrete mathematics. Experienced in Object-Oriented design, analysis, and prototyping in teams of 3-6 developers.\"\n$r = $d.Range($summary.Range.Start, $summary.Range.End - 1)\n$r.Text = $newSummary\n\n# --- 2. Remove the trailing \"Excel Forecasting Simulations\" table row ----\n$t = $d.Tables.Item(1)\n$t.Rows.Item($t.Rows.Count).Delete()\n\n# --- 3. Re-even the remaining column widths (3116 -> 3117 twips) ---------\n$t.Columns.Item(1).Width = 155.85\n"}
